$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the PAT (Personal Access Token) values out of column D, rows 2-5.
$ws.Range("D2:D5").ClearContents()

# Row heights collapse back to the (single-line) default now that the long
# wrapped PAT text is gone.
$ws.Range("A2:E5").RowHeight = 13

# Move the active selection.
$ws.Range("C10").Select() | Out-Null
